$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 7
    4  = 3
    5  = 3
    6  = 5
    7  = 1
    8  = 4
    9  = 2
    10 = 2
    11 = 5
    12 = 8
    13 = 4
    14 = 3
    15 = 3
    16 = 4
    17 = 0
    18 = 3
    19 = 2
    20 = 2
    21 = 3
    22 = 8
    23 = 3
    24 = 3
    25 = 2
    26 = 1
    27 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
